$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 0.01461682877891991
$ws.Cells.Item(2, 2).Value = -0.6536770345289049
$ws.Cells.Item(2, 3).Value = 0.01460360101157396
$ws.Cells.Item(2, 4).Value = 2.90410644732146
$ws.Cells.Item(2, 5).Value = -5.450237561760328

$ws.Cells.Item(3, 1).Value = 0.1728643399999989
$ws.Cells.Item(3, 2).Value = 1.726891795775603
$ws.Cells.Item(3, 3).Value = 0.05345242601967006
$ws.Cells.Item(3, 4).Value = 232.5090162211441
$ws.Cells.Item(3, 5).Value = -341.2981389385268

$ws.Cells.Item(4, 1).Value = 0.1686058266666655
$ws.Cells.Item(4, 2).Value = 1.106598544986156
$ws.Cells.Item(4, 3).Value = 0.05902422017440354
$ws.Cells.Item(4, 4).Value = 284.8276825796112
$ws.Cells.Item(4, 5).Value = -221.5776036302574

$ws.Cells.Item(5, 1).Value = 0.1524477599999987
$ws.Cells.Item(5, 2).Value = 0.5709679500941993
$ws.Cells.Item(5, 3).Value = 0.05109702746579042
$ws.Cells.Item(5, 4).Value = 280.0321931813189
$ws.Cells.Item(5, 5).Value = -266.9804431420322

$ws.Cells.Item(6, 1).Value = 0.1166450242424227
$ws.Cells.Item(6, 2).Value = 0.8326149783915795
$ws.Cells.Item(6, 3).Value = 0.03409674077543669
$ws.Cells.Item(6, 4).Value = 340.8115996962164
$ws.Cells.Item(6, 5).Value = -308.2994365764367

$ws.Cells.Item(7, 1).Value = 0.1310955818181801
$ws.Cells.Item(7, 2).Value = 1.045055280971326
$ws.Cells.Item(7, 3).Value = 0.03513250032469947
$ws.Cells.Item(7, 4).Value = 340.3935972893077
$ws.Cells.Item(7, 5).Value = -299.8507895223406

$ws.Cells.Item(8, 1).Value = 0.1154471466666654
$ws.Cells.Item(8, 2).Value = 0.5764486838405825
$ws.Cells.Item(8, 3).Value = 0.03401543262826115
$ws.Cells.Item(8, 4).Value = 397.9325953967981
$ws.Cells.Item(8, 5).Value = -259.4258982479697

$ws.Cells.Item(9, 1).Value = 0.0755528266666657
$ws.Cells.Item(9, 2).Value = 0.1054074822819475
$ws.Cells.Item(9, 3).Value = 0.02351459052264597
$ws.Cells.Item(9, 4).Value = 192.1938709167994
$ws.Cells.Item(9, 5).Value = -176.803078937943

$ws.Cells.Item(10, 1).Value = 0.08555979666666559
$ws.Cells.Item(10, 2).Value = 0.4343685029898836
$ws.Cells.Item(10, 3).Value = 0.02585716704468571
$ws.Cells.Item(10, 4).Value = 192.3915993930445
$ws.Cells.Item(10, 5).Value = -186.142059491945

$ws.Cells.Item(11, 1).Value = 0.06898184333333222
$ws.Cells.Item(11, 2).Value = 0.1067754478097014
$ws.Cells.Item(11, 3).Value = 0.02190872016322986
$ws.Cells.Item(11, 4).Value = 197.2585413971364
$ws.Cells.Item(11, 5).Value = -192.34991899551

$ws.Cells.Item(12, 1).Value = 0.05163072666666586
$ws.Cells.Item(12, 2).Value = -0.411430301099152
$ws.Cells.Item(12, 3).Value = 0.0215885898974164
$ws.Cells.Item(12, 4).Value = 91.81050562984844
$ws.Cells.Item(12, 5).Value = -128.9796578485333

$ws.Cells.Item(13, 1).Value = 0.04136679666666598
$ws.Cells.Item(13, 2).Value = -0.1045325132266157
$ws.Cells.Item(13, 3).Value = 0.01477798672358223
$ws.Cells.Item(13, 4).Value = 104.408713663341
$ws.Cells.Item(13, 5).Value = -115.2323676483018

$ws.Cells.Item(14, 1).Value = 0.05533729666666595
$ws.Cells.Item(14, 2).Value = -0.4143575405123757
$ws.Cells.Item(14, 3).Value = 0.02240326922942375
$ws.Cells.Item(14, 4).Value = 97.00017227548516
$ws.Cells.Item(14, 5).Value = -120.6684954716332

$ws.Cells.Item(15, 1).Value = 0.03083268541666636
$ws.Cells.Item(15, 2).Value = -0.3530350083615316
$ws.Cells.Item(15, 3).Value = 0.01627834077565098
$ws.Cells.Item(15, 4).Value = 67.85116839020948
$ws.Cells.Item(15, 5).Value = -82.37334929515534

$ws.Cells.Item(16, 1).Value = 0.02983729999999946
$ws.Cells.Item(16, 2).Value = -0.2395837204645644
$ws.Cells.Item(16, 3).Value = 0.01230040816995616
$ws.Cells.Item(16, 4).Value = 72.42107978297658
$ws.Cells.Item(16, 5).Value = -74.35767362446802

$ws.Cells.Item(17, 1).Value = 0.03187929999999949
$ws.Cells.Item(17, 2).Value = -0.2900484156769752
$ws.Cells.Item(17, 3).Value = 0.01450449659839718
$ws.Cells.Item(17, 4).Value = 61.9403494263091
$ws.Cells.Item(17, 5).Value = -81.33530838608252

$ws.Cells.Item(18, 1).Value = 0.02738696666666629
$ws.Cells.Item(18, 2).Value = -0.4952661885593346
$ws.Cells.Item(18, 3).Value = 0.01410683746600191
$ws.Cells.Item(18, 4).Value = 21.88413086905823
$ws.Cells.Item(18, 5).Value = -73.64613924139852

$ws.Cells.Item(19, 1).Value = 0.02356666666666625
$ws.Cells.Item(19, 2).Value = -0.4831209571760536
$ws.Cells.Item(19, 3).Value = 0.01344244939790623
$ws.Cells.Item(19, 4).Value = 29.93576416041296
$ws.Cells.Item(19, 5).Value = -66.30320178915427

$ws.Cells.Item(20, 1).Value = 0.0176335428571426
$ws.Cells.Item(20, 2).Value = -0.4193371263871642
$ws.Cells.Item(20, 3).Value = 0.01066363981611673
$ws.Cells.Item(20, 4).Value = 31.99562501018344
$ws.Cells.Item(20, 5).Value = -76.24182983192375

$ws.Cells.Item(21, 1).Value = 0.01688125520833304
$ws.Cells.Item(21, 2).Value = -0.4759898255524047
$ws.Cells.Item(21, 3).Value = 0.01193641439824343
$ws.Cells.Item(21, 4).Value = 28.25464481318303
$ws.Cells.Item(21, 5).Value = -84.48462977081785

$ws.Cells.Item(22, 1).Value = 0.0256956999999996
$ws.Cells.Item(22, 2).Value = -0.516632783353136
$ws.Cells.Item(22, 3).Value = 0.0143637551531884
$ws.Cells.Item(22, 4).Value = 38.00209543046216
$ws.Cells.Item(22, 5).Value = -61.93016930583573

$ws.Cells.Item(23, 1).Value = 0.02085103854166638
$ws.Cells.Item(23, 2).Value = -0.4791163098488145
$ws.Cells.Item(23, 3).Value = 0.01254157844963881
$ws.Cells.Item(23, 4).Value = 21.1960434028457
$ws.Cells.Item(23, 5).Value = -72.19888026138942

$ws.Cells.Item(24, 1).Value = 0.0198869323601969
$ws.Cells.Item(24, 2).Value = -0.6493123890965078
$ws.Cells.Item(24, 3).Value = 0.01457858597115233
$ws.Cells.Item(24, 4).Value = 38.06348417639627
$ws.Cells.Item(24, 5).Value = -56.92832554438967

$ws.Cells.Item(25, 1).Value = 0.0201829787709797
$ws.Cells.Item(25, 2).Value = -0.667705284466882
$ws.Cells.Item(25, 3).Value = 0.01529077329480307
$ws.Cells.Item(25, 4).Value = 29.37582716036991
$ws.Cells.Item(25, 5).Value = -59.75109517053831

$ws.Cells.Item(26, 1).Value = 0.01942165245889042
$ws.Cells.Item(26, 2).Value = -0.6041265854581094
$ws.Cells.Item(26, 3).Value = 0.01393560204946034
$ws.Cells.Item(26, 4).Value = 46.48894368597867
$ws.Cells.Item(26, 5).Value = -67.5379172464626

$ws.Cells.Item(27, 1).Value = 0.01717669349123039
$ws.Cells.Item(27, 2).Value = -0.650699980811523
$ws.Cells.Item(27, 3).Value = 0.01486336823295297
$ws.Cells.Item(27, 4).Value = 33.33084147069215
$ws.Cells.Item(27, 5).Value = -68.31909724871178

$ws.Cells.Item(28, 1).Value = 0.02122841387036161
$ws.Cells.Item(28, 2).Value = -0.6803653909792889
$ws.Cells.Item(28, 3).Value = 0.01608332701901606
$ws.Cells.Item(28, 4).Value = 37.76875511559213
$ws.Cells.Item(28, 5).Value = -51.75771851419073

$ws.Cells.Item(29, 1).Value = 0.01304621389018141
$ws.Cells.Item(29, 2).Value = -0.460428630672934
$ws.Cells.Item(29, 3).Value = 0.01086760150764514
$ws.Cells.Item(29, 4).Value = 36.81280712484585
$ws.Cells.Item(29, 5).Value = -55.48258067246545

$ws.Cells.Item(30, 1).Value = 0.0171354159736838
$ws.Cells.Item(30, 2).Value = -0.557945038109792
$ws.Cells.Item(30, 3).Value = 0.01327944039906333
$ws.Cells.Item(30, 4).Value = 32.62328107762917
$ws.Cells.Item(30, 5).Value = -65.38925204724993

$ws.Cells.Item(31, 1).Value = 0.01669786562914701
$ws.Cells.Item(31, 2).Value = -0.5791854679963012
$ws.Cells.Item(31, 3).Value = 0.01377611248523513
$ws.Cells.Item(31, 4).Value = 30.12708172981849
$ws.Cells.Item(31, 5).Value = -70.00280588147028

$ws.Cells.Item(32, 1).Value = 0.02460342249999979
$ws.Cells.Item(32, 2).Value = -1.00242956574095
$ws.Cells.Item(32, 3).Value = 0.02196530719355607
$ws.Cells.Item(32, 4).Value = 8.958492269745213
$ws.Cells.Item(32, 5).Value = -12.1072120678675

$ws.Cells.Item(33, 1).Value = 0.01531722854212956
$ws.Cells.Item(33, 2).Value = -0.7023149468219593
$ws.Cells.Item(33, 3).Value = 0.01645038932621357
$ws.Cells.Item(33, 4).Value = 17.33980085159854
$ws.Cells.Item(33, 5).Value = -92.01874341558198

$ws.Cells.Item(34, 1).Value = 0.02005938349849735
$ws.Cells.Item(34, 2).Value = -0.8371256729959542
$ws.Cells.Item(34, 3).Value = 0.01795984244250424
$ws.Cells.Item(34, 4).Value = 19.53146198129168
$ws.Cells.Item(34, 5).Value = -86.86632210679339

$ws.Cells.Item(35, 1).Value = 0.01970191438883411
$ws.Cells.Item(35, 2).Value = -0.7928128300649632
$ws.Cells.Item(35, 3).Value = 0.01857757543702953
$ws.Cells.Item(35, 4).Value = 13.67992466406772
$ws.Cells.Item(35, 5).Value = -108.7760421171383

$ws.Cells.Item(36, 1).Value = 0.01616103368388381
$ws.Cells.Item(36, 2).Value = -0.7344788531821792
$ws.Cells.Item(36, 3).Value = 0.01680391572955407
$ws.Cells.Item(36, 4).Value = 19.44638498139446
$ws.Cells.Item(36, 5).Value = -82.5016456698483

$ws.Cells.Item(37, 1).Value = 0.01818225848422167
$ws.Cells.Item(37, 2).Value = -0.7876842709308629
$ws.Cells.Item(37, 3).Value = 0.01713919452276819
$ws.Cells.Item(37, 4).Value = 11.37741414608985
$ws.Cells.Item(37, 5).Value = -90.38536625301617

$ws.Cells.Item(38, 1).Value = 0.01969223741332427
$ws.Cells.Item(38, 2).Value = -0.8308961020123464
$ws.Cells.Item(38, 3).Value = 0.01939952079340379
$ws.Cells.Item(38, 4).Value = 13.41435035713817
$ws.Cells.Item(38, 5).Value = -101.2166942433644

$ws.Cells.Item(39, 1).Value = 0.02196707833790037
$ws.Cells.Item(39, 2).Value = -0.8473047673472127
$ws.Cells.Item(39, 3).Value = 0.01817307655378765
$ws.Cells.Item(39, 4).Value = 16.28047606733938
$ws.Cells.Item(39, 5).Value = -82.83675410526831

$ws.Cells.Item(40, 1).Value = 0.02534582564787492
$ws.Cells.Item(40, 2).Value = -1.087081346340354
$ws.Cells.Item(40, 3).Value = 0.0242421500529274
$ws.Cells.Item(40, 4).Value = 13.65391634351252
$ws.Cells.Item(40, 5).Value = -7.120863890885293

$ws.Cells.Item(41, 1).Value = 0.01331892569291484
$ws.Cells.Item(41, 2).Value = -0.5294083889143651
$ws.Cells.Item(41, 3).Value = 0.01218165047656733
$ws.Cells.Item(41, 4).Value = 18.04556296199106
$ws.Cells.Item(41, 5).Value = -72.9745672578326

$ws.Cells.Item(42, 1).Value = 0.01866133842757729
$ws.Cells.Item(42, 2).Value = -0.7226529194067505
$ws.Cells.Item(42, 3).Value = 0.01705057858899307
$ws.Cells.Item(42, 4).Value = 18.010380305741
$ws.Cells.Item(42, 5).Value = -88.34238076227761

